$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A47').Value = 8.0
$ws.Range('B47').Value = 'WHITESPACE/NEW_LINE'
$ws.Range('C47').Value = '\n'
$ws.Range('A48').Value = 8.0
$ws.Range('B48').Value = 'WHITESPACE/SPACE'
$ws.Range('C48').Value = '空格'
$ws.Range('A49').Value = 8.0
$ws.Range('B49').Value = 'WHITESPACE/SPACE'
$ws.Range('C49').Value = '空格'
$ws.Range('A50').Value = 8.0
$ws.Range('B50').Value = 'WHITESPACE/SPACE'
$ws.Range('C50').Value = '空格'
$ws.Range('A51').Value = 8.0
$ws.Range('B51').Value = 'WHITESPACE/SPACE'
$ws.Range('C51').Value = '空格'
$ws.Range('A52').Value = 8.0
$ws.Range('B52').Value = 'IDENTIFIER'
$ws.Range('C52').Value = 'v1'
$ws.Range('A53').Value = 8.0
$ws.Range('B53').Value = 'WHITESPACE/SPACE'
$ws.Range('C53').Value = '空格'
$ws.Range('A54').Value = 8.0
$ws.Range('B54').Value = 'DOUBLECHARDELIMITER/ASSIGN'
$ws.Range('C54').Value = ':='
$ws.Range('A55').Value = 8.0
$ws.Range('B55').Value = 'WHITESPACE/SPACE'
$ws.Range('C55').Value = '空格'
$ws.Range('A56').Value = 8.0
$ws.Range('B56').Value = 'IDENTIFIER'
$ws.Range('C56').Value = 'v1'
$ws.Range('A57').Value = 8.0
$ws.Range('B57').Value = 'WHITESPACE/SPACE'
$ws.Range('C57').Value = '空格'
$ws.Range('A58').Value = 8.0
$ws.Range('B58').Value = 'OPERATOR/PLUS'
$ws.Range('C58').Value = '+'
$ws.Range('A59').Value = 8.0
$ws.Range('B59').Value = 'WHITESPACE/SPACE'
$ws.Range('C59').Value = '空格'
$ws.Range('A60').Value = 8.0
$ws.Range('B60').Value = 'CONSTANT/UNSIGNED_INTEGER'
$ws.Range('C60').Value = '10'
$ws.Range('A61').Value = 8.0
$ws.Range('B61').Value = 'ENDDELIMITER/SEMICOLON'
$ws.Range('C61').Value = ';'
$ws.Range('A62').Value = 9.0
$ws.Range('B62').Value = 'WHITESPACE/NEW_LINE'
$ws.Range('C62').Value = '\n'
$ws.Range('A63').Value = 9.0
$ws.Range('B63').Value = 'WHITESPACE/SPACE'
$ws.Range('C63').Value = '空格'
$ws.Range('A64').Value = 9.0
$ws.Range('B64').Value = 'WHITESPACE/SPACE'
$ws.Range('C64').Value = '空格'
$ws.Range('A65').Value = 9.0
$ws.Range('B65').Value = 'WHITESPACE/SPACE'
$ws.Range('C65').Value = '空格'
$ws.Range('A66').Value = 9.0
$ws.Range('B66').Value = 'WHITESPACE/SPACE'
$ws.Range('C66').Value = '空格'
$ws.Range('A67').Value = 9.0
$ws.Range('B67').Value = 'COMMENT'
$ws.Range('C67').Value = '{write(''a'''')}'
$ws.Range('A68').Value = 10.0
$ws.Range('B68').Value = 'WHITESPACE/NEW_LINE'
$ws.Range('C68').Value = '\n'
$ws.Range('A69').Value = 10.0
$ws.Range('B69').Value = 'WHITESPACE/SPACE'
$ws.Range('C69').Value = '空格'
$ws.Range('A70').Value = 10.0
$ws.Range('B70').Value = 'WHITESPACE/SPACE'
$ws.Range('C70').Value = '空格'
$ws.Range('A71').Value = 10.0
$ws.Range('B71').Value = 'WHITESPACE/SPACE'
$ws.Range('C71').Value = '空格'
$ws.Range('A72').Value = 10.0
$ws.Range('B72').Value = 'WHITESPACE/SPACE'
$ws.Range('C72').Value = '空格'
$ws.Range('A73').Value = 10.0
$ws.Range('B73').Value = 'WORD/WRITE'
$ws.Range('C73').Value = 'write'
$ws.Range('A74').Value = 10.0
$ws.Range('B74').Value = 'PAIRDELIMITER/L_PARENTHESIS'
$ws.Range('C74').Value = '('
$ws.Range('A75').Value = 10.0
$ws.Range('B75').Value = 'IDENTIFIER'
$ws.Range('C75').Value = 'v1'
$ws.Range('A76').Value = 10.0
$ws.Range('B76').Value = 'PAIRDELIMITER/R_PARENTHESIS'
$ws.Range('C76').Value = ')'
$ws.Range('A77').Value = 10.0
$ws.Range('B77').Value = 'ENDDELIMITER/SEMICOLON'
$ws.Range('C77').Value = ';'
$ws.Range('A78').Value = 11.0
$ws.Range('B78').Value = 'WHITESPACE/NEW_LINE'
$ws.Range('C78').Value = '\n'
$ws.Range('A79').Value = 11.0
$ws.Range('B79').Value = 'WORD/END'
$ws.Range('C79').Value = 'end'
$ws.Range('A80').Value = 11.0
$ws.Range('B80').Value = 'PROGRAMEND'
$ws.Range('C80').Value = '.'

$ws.Rows.Item(81).Delete()
